# Jogos_da_Semana_FlashScore_2024-11-26.xlsx update
# - Remove the BRAZIL - SERIE A BETANO "Fortaleza vs Flamengo RJ" match (row 6),
#   which shifts the following rows (Atletico-MG/Juventude, Palmeiras/Botafogo RJ,
#   Santa Fe/Millonarios, Progreso/Penarol) up by one row.
# - Refresh a handful of odds values on the remaining early rows (2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 6 (entire row, cells shift up)
$ws.Rows(6).Delete()

# --- Row 2 (Sarmiento Junin vs Platense) odds refresh ---
$ws.Range("O2").Value = 1.8
$ws.Range("P2").Value = 1.91
$ws.Range("Q2").Value = 3.6
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 2

# --- Row 3 (Union de Santa Fe vs Talleres Cordoba) odds refresh ---
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5

# --- Row 4 (Guabira vs Tomayapo) odds refresh ---
$ws.Range("G4").Value = 1.7
$ws.Range("I4").Value = 4.5
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 9
$ws.Range("AH4").Value = 23
$ws.Range("AL4").Value = 34
$ws.Range("AO4").Value = 17
$ws.Range("AY4").Value = 81
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 101
$ws.Range("BD4").Value = 151

# --- Row 5 (Blooming vs Universitario de Vinto) odds refresh ---
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.2
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.08
$ws.Range("W5").Value = 8
$ws.Range("Z5").Value = 13
$ws.Range("AD5").Value = 7.5
$ws.Range("AG5").Value = 15
$ws.Range("BD5").Value = 201
